# Automatische test-sync: 2025-07-29 22:06:50
# Append the Testmail #18 log entry to the "Logs" sheet and refresh the
# "Dashboard" category counts to match.

$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")
$dash = $wb.Worksheets.Item("Dashboard")

# New row of data for the Logs sheet (row 20)
$row = 20
$logs.Cells.Item($row, 1).Value = "Bestel je 200 stuks M8-bouten RVS voor Van Dijk?"
$logs.Cells.Item($row, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item($row, 3).Value = "Testmail #18: Bestel je 200 stuks M8-bouten RVS voor Van Dijk?"
$logs.Cells.Item($row, 4).Value = "Bestelling / Levering"
$logs.Cells.Item($row, 5).Value = "Beste klant,`nBedankt voor uw e-mail. Het spijt me, maar het lijkt erop dat u per ongeluk een testmail heeft gestuurd. Als u daadwerkelijk 200 stuks M8-bouten RVS wilt bestellen voor Van Dijk, raad ik u aan om contact op te nemen met onze verkoopafdeling via [verkoop@email.com] of telefonisch via [telefoonnummer]. Zij helpen u graag verder met uw bestelling.`nMet vriendelijke groet,`n[Naam] Nederlandse e-mailassistent van <bedrijfsnaam>"
$logs.Cells.Item($row, 6).Value = "2025-07-29 22:06:31"
$logs.Cells.Item($row, 7).Value = "Ja"
$logs.Cells.Item($row, 8).Value = "Nee"
$logs.Cells.Item($row, 9).Value = "Ja"
$logs.Cells.Item($row, 10).Value = "Nee"

# Update the Dashboard summary table: the new entry is another
# "Bestelling / Levering" mail, so that category now counts 2, matching
# (and swapping display order with) "Retour / Terugbetaling".
$dash.Cells.Item(5, 1).Value = "Bestelling / Levering"
$dash.Cells.Item(5, 2).Value = 2
$dash.Cells.Item(6, 1).Value = "Retour / Terugbetaling"
$dash.Cells.Item(6, 2).Value = 2

# Extend the conditional formatting sqref on every flagged column
# (D, G, H, I, J) from row 19 down to the newly added row 20.
foreach ($col in @("D","G","H","I","J")) {
    $oldRange = $logs.Range($col + "2:" + $col + "19")
    $newRange = $logs.Range($col + "2:" + $col + "20")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}
